$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The "Price" column (D) stores numeric-looking values as plain text
# (e.g. "318.89"). Force Text number format on each updated Price cell
# right before writing its new value so Excel does not auto-convert the
# string into a numeric value.

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "42.580.94"
$ws.Range("E2").Value = "  -2.21%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.336.16"
$ws.Range("E3").Value = "  -3.20%  "
$ws.Range("E4").Value = "  -0.06%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "318.89"
$ws.Range("E5").Value = "  -1.08%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "105.07"
$ws.Range("E6").Value = "  +1.07%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.636"
$ws.Range("E7").Value = "  -1.36%  "
$ws.Range("E8").Value = "  +0.06%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.613"
$ws.Range("E9").Value = "  -6.03%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "40.53"
$ws.Range("E10").Value = "  -3.20%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0919"
$ws.Range("E11").Value = "  -2.86%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "8.31"
$ws.Range("E12").Value = "  -4.37%  "
$ws.Range("E13").Value = "  -0.68%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.984"
$ws.Range("E14").Value = "  -4.85%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "15.86"
$ws.Range("E15").Value = "  -8.86%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "2.688.54"
$ws.Range("E16").Value = "  -3.26%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "2.322.75"
$ws.Range("E17").Value = "  -3.44%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "42.551.95"
$ws.Range("E18").Value = "  -2.30%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "7.67"
$ws.Range("E19").Value = "  +3.79%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.0000106"
$ws.Range("E20").Value = "  -4.10%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "76.52"
$ws.Range("E21").Value = "  +1.43%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "3.55"
$ws.Range("E22").Value = "  +2.37%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "258.55"
$ws.Range("E23").Value = "  -0.82%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.31"
$ws.Range("E24").Value = "  -5.38%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "9.67"
$ws.Range("E25").Value = "  -0.55%  "
$ws.Range("E26").Value = "  +0.12%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "11.33"
$ws.Range("E27").Value = "  -5.64%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "22.97"
$ws.Range("E28").Value = "  -0.03%  "
$ws.Range("E29").Value = "  +1.82%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "174.47"
$ws.Range("E30").Value = "  -2.50%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "35.68"
$ws.Range("E31").Value = "  -6.53%  "
$ws.Range("B32").Value = "Hedera"
$ws.Range("C32").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.0889"
$ws.Range("E32").Value = "  -5.07%  "
$ws.Range("B33").Value = "WEMIXToken"
$ws.Range("C33").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "2.98"
$ws.Range("E33").Value = "  -8.01%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "6.04"
$ws.Range("E34").Value = "  +0.92%  "
$ws.Range("E35").Value = "  -1.77%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.112"
$ws.Range("E36").Value = "  +5.57%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "4.52"
$ws.Range("E37").Value = "  -7.61%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.0353"
$ws.Range("E38").Value = "  -4.58%  "
$ws.Range("E39").Value = "  -5.16%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "2.61"
$ws.Range("E40").Value = "  -9.96%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "1.46"
$ws.Range("E41").Value = "  -10.24%  "
$ws.Range("B42").Value = "MultiversX"
$ws.Range("C42").Value = "https://coinranking.com/coin/omwkOTglq+multiversx-egld"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "70.25"
$ws.Range("E42").Value = "  +1.55%  "
$ws.Range("B43").Value = "Algorand"
$ws.Range("C43").Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.231"
$ws.Range("E43").Value = "  -1.30%  "
$ws.Range("E44").Value = "  -0.14%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "115.03"
$ws.Range("E45").Value = "  -7.81%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "11.75"
$ws.Range("E46").Value = "  -7.17%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "5.49"
$ws.Range("E47").Value = "  -3.50%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "9.12"
$ws.Range("E48").Value = "  -4.73%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "84.56"
$ws.Range("E49").Value = "  +9.81%  "
$ws.Range("E50").Value = "  +1.97%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.0994"
$ws.Range("E51").Value = "  -1.68%  "
